# Insert a new week of "Chirimoya" price data (2021-11-04 / serial 44504) at
# the top of the Femacal de La Calera block, pushing the existing rows 84:99
# down to 87:102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 84 (existing rows 84:99 shift down to 87:102).
$ws.Range("A84:T86").EntireRow.Insert()

# New data rows (84:86) - same market/region/product metadata as the rest of
# the block, new date + volumes/prices + origin "Provincia de Limari".
$newRows = @(
    @{ Row = 84; Quality = "Especial"; Volume = 54; Min = 26000; Max = 26000; Avg = 26000; PerKg = 2600 },
    @{ Row = 85; Quality = "Primera";  Volume = 57; Min = 24000; Max = 24000; Avg = 24000; PerKg = 2400 },
    @{ Row = 86; Quality = "Segunda";  Volume = 50; Min = 20000; Max = 20000; Avg = 20000; PerKg = 2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0).Date
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Volume
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Avg
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.PerKg
    $ws.Cells.Item($row, 20).Value = 10
}
